$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.674.14'
$ws.Range('E2').Value = '  -2.95%  '
$ws.Range('D3').Value = '3.362.40'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.23'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.09'
$ws.Range('E6').Value = '  +7.74%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.362.06'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.476'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.57'
$ws.Range('E10').Value = '  +4.15%  '
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.389'
$ws.Range('E12').Value = '  +4.13%  '
$ws.Range('D13').Value = '3.933.20'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.365.00'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000172'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.13'
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('D18').Value = '60.874.13'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('E19').Value = '  +6.36%  '
$ws.Range('E20').Value = '  +3.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.44'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '371.60'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.567'
$ws.Range('E23').Value = '  +2.68%  '
$ws.Range('D24').Value = '3.496.39'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.57'
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('E27').Value = '  +10.87%  '
$ws.Range('E28').Value = '  +23.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.70'
$ws.Range('E29').Value = '  +11.86%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  +4.47%  '
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('E33').Value = '  +4.12%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '3.392.62'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.36'
$ws.Range('E36').Value = '  +3.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.53'
$ws.Range('E37').Value = '  +4.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.91'
$ws.Range('E38').Value = '  +4.63%  '
$ws.Range('E39').Value = '  +5.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '163.00'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('E41').Value = '  +4.38%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.39'
$ws.Range('E43').Value = '  +4.16%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.20'
$ws.Range('E44').Value = '  +11.75%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.758'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.24'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('E47').Value = '  +3.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.12'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.97'
$ws.Range('E49').Value = '  +5.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.16'
$ws.Range('E50').Value = '  +15.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.43'
$ws.Range('E51').Value = '  +14.20%  '
